{"js": "// Replace the four Japanese sentences with their Portuguese translations.\nconst replacements = [\n  {\n    find: \"\u3053\u306e\u53e4\u304f\u304b\u3089\u4f1d\u308f\u308b\u30ec\u30b7\u30d4\u306f\u3001\u4e00\u53e3\u98f2\u3080\u3054\u3068\u306b\u672c\u683c\u7684\u3067\u3057\u3063\u304b\u308a\u3068\u3057\u305f\u5473\u308f\u3044\u3092\u7d04\u675f\u3057\u307e\u3059\u3002\",\n    replace: \"Esta receita centen\u00e1ria promete um sabor aut\u00eantico e robusto em cada gole.\"\n  },\n  {\n    find: \"\u98a8\u5473\u306f\u5f37\u70c8\u3067\u3042\u308a\u306a\u304c\u3089\u30d0\u30e9\u30f3\u30b9\u304c\u53d6\u308c\u3066\u304a\u308a\u3001\u5feb\u9069\u3067\u5fc3\u5730\u3088\u3044\u4f53\u9a13\u3092\u751f\u307f\u51fa\u3057\u307e\u3059\u3002\",\n    replace: \"Os sabores s\u00e3o intensos, mas equilibrados, proporcionando uma experi\u00eancia reconfortante e relaxante.\"\n  },\n  {\n    find: \"\u304a\u597d\u307f\u306e\u65b9\u6cd5\u3067\u30c1\u30e3\u30a4\u3092\u304a\u697d\u3057\u307f\u3044\u305f\u3060\u3051\u308b\u3088\u3046\u3001\u7c21\u5358\u306a\u6df9\u308c\u65b9\u306e\u8aac\u660e\u66f8\u304c\u4ed8\u5c5e\u3057\u3066\u3044\u307e\u3059\u3002\",\n    replace: \"Instru\u00e7\u00f5es simples de preparo est\u00e3o inclu\u00eddas para ajud\u00e1-lo a saborear seu chai exatamente do jeito que voc\u00ea gosta.\"\n  },\n  {\n    find: \"Mystic Spice Chai Tea \u304c\u304a\u5ba2\u69d8\u306e\u3054\u671f\u5f85\u306b\u6dfb\u3048\u306a\u3044\u5834\u5408\u306f\u3001\u5f53\u793e\u304c\u6539\u5584\u3059\u308b\u3088\u3046\u52aa\u3081\u307e\u3059\u3002\",\n    replace: \"Se o Mystic Spice Chai Tea n\u00e3o atender suas expectativas, estamos comprometidos em resolver da melhor maneira poss\u00edvel.\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the four Japanese sentences with their Portuguese translations.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Find    = \"\u3053\u306e\u53e4\u304f\u304b\u3089\u4f1d\u308f\u308b\u30ec\u30b7\u30d4\u306f\u3001\u4e00\u53e3\u98f2\u3080\u3054\u3068\u306b\u672c\u683c\u7684\u3067\u3057\u3063\u304b\u308a\u3068\u3057\u305f\u5473\u308f\u3044\u3092\u7d04\u675f\u3057\u307e\u3059\u3002\"\n        Replace = \"Esta receita centen\u00e1ria promete um sabor aut\u00eantico e robusto em cada gole.\"\n    },\n    @{\n        Find    = \"\u98a8\u5473\u306f\u5f37\u70c8\u3067\u3042\u308a\u306a\u304c\u3089\u30d0\u30e9\u30f3\u30b9\u304c\u53d6\u308c\u3066\u304a\u308a\u3001\u5feb\u9069\u3067\u5fc3\u5730\u3088\u3044\u4f53\u9a13\u3092\u751f\u307f\u51fa\u3057\u307e\u3059\u3002\"\n        Replace = \"Os sabores s\u00e3o intensos, mas equilibrados, proporcionando uma experi\u00eancia reconfortante e relaxante.\"\n    },\n    @{\n        Find    = \"\u304a\u597d\u307f\u306e\u65b9\u6cd5\u3067\u30c1\u30e3\u30a4\u3092\u304a\u697d\u3057\u307f\u3044\u305f\u3060\u3051\u308b\u3088\u3046\u3001\u7c21\u5358\u306a\u6df9\u308c\u65b9\u306e\u8aac\u660e\u66f8\u304c\u4ed8\u5c5e\u3057\u3066\u3044\u307e\u3059\u3002\"\n        Replace = \"Instru\u00e7\u00f5es simples de preparo est\u00e3o inclu\u00eddas para ajud\u00e1-lo a saborear seu chai exatamente do jeito que voc\u00ea gosta.\"\n    },\n    @{\n        Find    = \"Mystic Spice Chai Tea \u304c\u304a\u5ba2\u69d8\u306e\u3054\u671f\u5f85\u306b\u6dfb\u3048\u306a\u3044\u5834\u5408\u306f\u3001\u5f53\u793e\u304c\u6539\u5584\u3059\u308b\u3088\u3046\u52aa\u3081\u307e\u3059\u3002\"\n        Replace = \"Se o Mystic Spice Chai Tea n\u00e3o atender suas expectativas, estamos comprometidos em resolver da melhor maneira poss\u00edvel.\"\n    }\n)\n\nforeach ($item in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $item.Find\n    $find.Replacement.Text = $item.Replace\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
